$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7").Value = -7.472999999999999
$ws.Range("C8").Value = -12.542
$ws.Range("C10").Value = -12.889
$ws.Range("C12").Value = -11.207
$ws.Range("D14").Value = -7.878
$ws.Range("D15").Value = -8.317000000000002
$ws.Range("C18").Value = -14.011
$ws.Range("D18").Value = -8.431000000000001
$ws.Range("D20").Value = -7.19
$ws.Range("C25").Value = -11.9
$ws.Range("D29").Value = -7.292
$ws.Range("D30").Value = -7.231
$ws.Range("D31").Value = -7.938000000000001
$ws.Range("D35").Value = -7.582000000000001
$ws.Range("C37").Value = -13.345
$ws.Range("D40").Value = -7.293000000000001
$ws.Range("D44").Value = -7.181
$ws.Range("D50").Value = -8.105
$ws.Range("D54").Value = -8.089
$ws.Range("C55").Value = -14.117
$ws.Range("C68").Value = -11.167
$ws.Range("D68").Value = -6.879
$ws.Range("D76").Value = -7.161
$ws.Range("C77").Value = -13.401
$ws.Range("C78").Value = -13.375
$ws.Range("C79").Value = -13.099
$ws.Range("C80").Value = -13.268
$ws.Range("C81").Value = -13.094
$ws.Range("C82").Value = -11.933
$ws.Range("C84").Value = -13.181
$ws.Range("D87").Value = -8.341000000000001
$ws.Range("D88").Value = -8.004999999999999
$ws.Range("D92").Value = -6.584000000000001
$ws.Range("D96").Value = -7.267
$ws.Range("D98").Value = -8.198000000000002
$ws.Range("C101").Value = -13.048
$ws.Range("D101").Value = -7.997
$ws.Range("C102").Value = -13.743
$ws.Range("D102").Value = -7.629

$wb.Save()
